# STATS.xlsx — "add extract all darfeuille"
#
# Fills in the second-pass ("Seconde passe") figures that were still
# missing for the "METRO ABC,AB,A,B,C,D COMMUNS" line (row 13), records
# an unparsable raw reading "3830?" for "RHONE EXPRESS COMMUNS" (row 21,
# which breaks its time-of-day formula), normalizes the number format on
# I23 to match the rest of column I, clears the two stray formatted-but-
# empty cells on the "PRET TOTAL" row (L29/M29) together with the color
# scale that targeted M29, and leaves the selection where the author left
# it (H22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: "METRO ABC,AB,A,B,C,D COMMUNS" second-pass figures -----------
$ws.Range("H13").Value = 2856
$ws.Range("I13").Formula = "=H13/86400"
$ws.Range("I13").NumberFormat = "hh:\ mm:\ ss"
$ws.Range("J13").Value = 39723
$ws.Range("K13").Value = 4626
$ws.Range("L13").Formula = "=J13+K13"

# --- Row 21: "RHONE EXPRESS COMMUNS" raw time reading is unparsable -------
# (turns I21's H21/86400 formula into a #VALUE! error, same as Excel would)
$ws.Range("H21").Value = "3830?"

# --- Row 23: align I23's number format with the rest of column I ---------
$ws.Range("I23").NumberFormat = "hh:\ mm:\ ss"

# --- Row 29: drop the stray empty L29/M29 cells and their conditional ----
# formatting (a leftover color scale that no longer applies to anything)
$ws.Range("M29").FormatConditions.Delete()
$ws.Range("L29").Clear()
$ws.Range("M29").Clear()

# --- Restore the author's last selection ----------------------------------
$null = $ws.Range("H22").Select()
